$d = $word.ActiveDocument

$pairs = @(
    @{old="781÷8=97, 5"; new="524÷7=74, 6"},
    @{old="541÷4=135, 1"; new="923÷7=131, 6"},
    @{old="774÷9=86, 0"; new="333÷5=66, 3"},
    @{old="987÷8=123, 3"; new="428÷4=107, 0"},
    @{old="542÷2=271, 0"; new="398÷6=66, 2"},
    @{old="291÷9=32, 3"; new="152÷5=30, 2"},
    @{old="488÷6=81, 2"; new="766÷8=95, 6"},
    @{old="474÷9=52, 6"; new="856÷6=142, 4"},
    @{old="934÷9=103, 7"; new="884÷7=126, 2"},
    @{old="604÷7=86, 2"; new="238÷8=29, 6"},
    @{old="158÷3=52, 2"; new="588÷6=98, 0"},
    @{old="850÷4=212, 2"; new="474÷5=94, 4"},
    @{old="211÷9=23, 4"; new="153÷6=25, 3"},
    @{old="576÷4=144, 0"; new="842÷8=105, 2"},
    @{old="945÷8=118, 1"; new="680÷8=85, 0"},
    @{old="842÷9=93, 5"; new="564÷9=62, 6"},
    @{old="462÷4=115, 2"; new="951÷7=135, 6"},
    @{old="267÷3=89, 0"; new="778÷8=97, 2"},
    @{old="353÷7=50, 3"; new="846÷5=169, 1"},
    @{old="486÷9=54, 0"; new="623÷8=77, 7"},
    @{old="713÷4=178, 1"; new="861÷9=95, 6"},
    @{old="465÷5=93, 0"; new="466÷3=155, 1"},
    @{old="882÷2=441, 0"; new="413÷9=45, 8"},
    @{old="663÷5=132, 3"; new="345÷7=49, 2"},
    @{old="614÷3=204, 2"; new="736÷9=81, 7"}
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2)
}
